# Test-User angelegt für schnelles LogIn
# Adds a new "test" user row (personnel_id 10) to the Tabelle1 user list so
# that logging in during development/testing is quicker.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the new user record in row 11 ---------------------------------
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "test"
$ws.Range("C11").Value = "Test"
$ws.Range("D11").Value = "Test"
$ws.Range("E11").Value = "Teststraße 1"
$ws.Range("F11").Value = 77777
$ws.Range("G11").Value = "Teststadt"
$ws.Range("H11").Value = "test.test@example.com"
$ws.Range("I11").Value = "test"
$ws.Range("J11").Value = 1
$ws.Range("K11").Value = $false

# password column keeps the same "hidden" custom number format as the rows above
$ws.Range("I11").NumberFormat = $ws.Range("I2").NumberFormat

# turn the e-mail address into a mailto: hyperlink, like every other row
$null = $ws.Hyperlinks.Add($ws.Range("H11"), "mailto:test.test@example.com")

# --- Misc. UI state that Excel persisted when the edit was made -----------
$null = $ws.Range("K18").Select()
